# Update the "Historico" news-history worksheet with the newsbot's latest
# scrape: clear a stray empty cell on the last existing row, then append
# three new rows (47-49) of news-item data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historico")

# Row 46 previously had an explicit-but-empty Publicado_em_BRT cell (B46);
# the refreshed export no longer writes that placeholder, so remove it.
$ws.Range("B46").ClearContents()

# --- Row 47 -----------------------------------------------------------
$ws.Range("A47").Value = "05/01/2026 07:50:59"
$ws.Range("B47").Value = "05/01 07:41"
$ws.Range("C47").Value = "g1 > Política"
$ws.Range("D47").Value = "Brasil deve pedir a palavra em reunião do Conselho de Segurança da ONU sobre Venezuela nesta segunda"
$ws.Range("E47").Value = "https://g1.globo.com/politica/noticia/2026/01/05/brasil-deve-pedir-a-palavra-em-reuniao-do-conselho-de-seguranca-da-onu-sobre-venezuela-nesta-segunda.ghtml"
$ws.Range("F47").Value = "lula"
$ws.Range("G47").Value = "as Nações Unidas, Sérgio Danese siga na linha do pronunciamento do presidente Luiz Inácio Lula da Silva (PT). `nAinda no sábado, Lula disse que a ação militar norte-americana em solo ve"

# --- Row 48 -----------------------------------------------------------
$ws.Range("A48").Value = "05/01/2026 07:51:00"
$ws.Range("C48").Value = "VEJA"
$ws.Range("D48").Value = "A leitura da Comissão de Relações Exteriores sobre ataques à Venezuela"
$ws.Range("E48").Value = "https://veja.abril.com.br/coluna/radar/a-leitura-da-comissao-de-relacoes-exteriores-sobre-ataques-a-venezuela/"

# --- Row 49 -----------------------------------------------------------
$ws.Range("A49").Value = "05/01/2026 07:51:00"
# Publicado_em_BRT is present-but-blank for this item (a lone text-prefix
# marker, mirroring the feed's own empty-field convention) rather than a
# fully absent cell.
$ws.Range("B49").Value = "'"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "VEJA"
$ws.Range("D49").Value = "TV Globo começa o ano dando um belo sacode em ‘Três Graças’"
$ws.Range("E49").Value = "https://veja.abril.com.br/coluna/veja-gente/tv-globo-comeca-o-ano-dando-um-belo-sacode-em-tres-gracas/"
$ws.Range("F49").Value = "ldo"
$ws.Range("G49").Value = "Novela de Aguina<b>ldo</b> Silva terá reviravoltas no capítulo desta segunda-feira, 5"
